$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("N2").Value = 2.28
$ws.Range("O2").Value = 1.62
$ws.Range("P2").Value = 1.37
$ws.Range("N4").Value = 3.55
$ws.Range("P4").Value = 1.85
$ws.Range("X4").Value = 15.5
$ws.Range("Z4").Value = 970
$ws.Range("G5").Value = 5.4
$ws.Range("J5").Value = 2.98
$ws.Range("S5").Value = 5.3
$ws.Range("F6").Value = 3
$ws.Range("H6").Value = 2.7
$ws.Range("I6").Value = 2.9
$ws.Range("N6").Value = 2.62
$ws.Range("T6").Value = 2
$ws.Range("U6").Value = 1.82
$ws.Range("V6").Value = 1.52
$ws.Range("W6").Value = 1.44
$ws.Range("AJ6").Value = 65
$ws.Range("P7").Value = 1.74
$ws.Range("I8").Value = 1.39
$ws.Range("P8").Value = 1.85
$ws.Range("T8").Value = 2.28
$ws.Range("V8").Value = 3.55
$ws.Range("AE8").Value = 970
$ws.Range("AF8").Value = 160
$ws.Range("N9").Value = 1.1
$ws.Range("P9").Value = 2.42
$ws.Range("R9").Value = 1.62
$ws.Range("F10").Value = 1.99
$ws.Range("G10").Value = 2.26
$ws.Range("H10").Value = 3.7
$ws.Range("I10").Value = 4.8
$ws.Range("J10").Value = 3.4
$ws.Range("K10").Value = 4
$ws.Range("P10").Value = 1.88
$ws.Range("Q10").Value = 1.91
$ws.Range("S10").Value = 3.1
$ws.Range("V10").Value = 1.29
$ws.Range("W10").Value = 1.8
$ws.Range("AB10").Value = 11
$ws.Range("AC10").Value = 9.4
$ws.Range("AF10").Value = 15.5
$ws.Range("AG10").Value = 12.5
$ws.Range("H11").Value = 1.69
$ws.Range("J11").Value = 3.5
$ws.Range("S11").Value = 3.9
$ws.Range("T11").Value = 2.12
$ws.Range("U11").Value = 1.72
$ws.Range("W11").Value = 1.14
$ws.Range("AI11").Value = 55
$ws.Range("T12").Value = 1.78
$ws.Range("U12").Value = 2.08
$ws.Range("P13").Value = 1.86
$ws.Range("AO13").Value = 18.5
$ws.Range("F14").Value = 4.5
$ws.Range("H14").Value = 1.92
$ws.Range("I14").Value = 1.97
$ws.Range("V14").Value = 2.02
$ws.Range("F15").Value = 1.63
$ws.Range("G15").Value = 1.69
$ws.Range("F16").Value = 1.34
$ws.Range("G16").Value = 1.37
$ws.Range("H16").Value = 11
$ws.Range("J16").Value = 5.6
$ws.Range("W16").Value = 3.7
$ws.Range("AD16").Value = 40
$ws.Range("AE16").Value = 180
$ws.Range("AH16").Value = 29
$ws.Range("AO16").Value = 230
$ws.Range("G17").Value = 2.96
$ws.Range("J17").Value = 2.56
$ws.Range("T17").Value = 2.22
$ws.Range("U17").Value = 1.66
$ws.Range("AF17").Value = 970
$ws.Range("I18").Value = 1.58
$ws.Range("T18").Value = 1.8
$ws.Range("U18").Value = 2.1
$ws.Range("V18").Value = 2.72
$ws.Range("Y18").Value = 10.5
$ws.Range("Z18").Value = 10.5
$ws.Range("AA18").Value = 14.5
$ws.Range("AB18").Value = 29
$ws.Range("AD18").Value = 10.5
$ws.Range("AF18").Value = 1000
$ws.Range("AH18").Value = 22
$ws.Range("AJ18").Value = 210
$ws.Range("AO18").Value = 6.8
$ws.Range("N19").Value = 3.3
$ws.Range("F20").Value = 1.44
$ws.Range("G20").Value = 1.46
$ws.Range("H20").Value = 8.6
$ws.Range("J20").Value = 4.7
$ws.Range("K20").Value = 5.1
$ws.Range("R20").Value = 1.4
$ws.Range("U20").Value = 1.8
$ws.Range("W20").Value = 3.1
$ws.Range("X20").Value = 20
$ws.Range("AB20").Value = 7.8
$ws.Range("AK20").Value = 970
$ws.Range("AM20").Value = 190
$ws.Range("J21").Value = 2.78
$ws.Range("T21").Value = 2.1
$ws.Range("U21").Value = 1.75
$ws.Range("G22").Value = 2.34
$ws.Range("H22").Value = 3.55
$ws.Range("O22").Value = 1.41
$ws.Range("S22").Value = 4.3
$ws.Range("T22").Value = 1.87
$ws.Range("U22").Value = 1.9
$ws.Range("W22").Value = 1.74
$ws.Range("X22").Value = 13
$ws.Range("Y22").Value = 12.5
$ws.Range("AA22").Value = 95
$ws.Range("AI22").Value = 75
$ws.Range("G23").Value = 3.15
$ws.Range("H23").Value = 2.68
$ws.Range("I23").Value = 2.74
$ws.Range("W23").Value = 1.46
$ws.Range("AD23").Value = 12.5
$ws.Range("AE23").Value = 34
$ws.Range("J24").Value = 2.78
$ws.Range("V24").Value = 1.59
